# TIMES-NZ ScenTrade_TRADE_PARMS.xlsx
#
# The two ~TFM_INS blocks on sheet "INS" (rows 3-12 and rows 16-25) had their
# Attribute/Year/AllRegions/Pset_PN values swapped between the blocks
# (effectively the second block's parameters were re-run and landed back in
# the first block's rows, and vice versa). This script rewrites the C
# (Attribute), D (Year - only present for the NCAP_PASTI row), F
# (AllRegions) and H (Pset_PN) cells of both blocks to restore/apply the
# correct values, without touching any other cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Txt($ws, $addr, $val) {
    # Force the cell to store the literal text given, even when it looks
    # like a number (Excel would otherwise auto-convert "0.0025" etc. into
    # a numeric value, losing the exact textual representation used by the
    # TIMES-NZ tooling).
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}

# row, Attribute(C), Year(D or $null), AllRegions(F), Pset_PN(H)
$rows = @(
    @(3,  "VAROM",      $null,   "0",                    "TB_ELC*"),
    @(4,  "LIFE",        $null,  "100",                  "TB_ELC*"),
    @(5,  "CAP2ACT",     $null,  "31.536",                "TB_ELC*"),
    @(6,  "PEAK(CON)",   $null,  "1",                    "TB_ELC*"),
    @(7,  "AFA",         $null,  "1",                    "TB_ELC*"),
    @(8,  "EFF",         $null,  "0.99",                  "TB_ELC*"),
    @(9,  "NCAP_ILED",   $null,  "3",                    "TB_ELC*"),
    @(10, "NCAP_PASTI",  "2010", "1.6",                   "TB_ELC_*"),
    @(11, "INVCOST",     $null,  "480.0",                 "TB_ELC*"),
    @(12, "FIXOM",       $null,  "0.96",                  "TB_ELC*"),
    @(16, "ACT_COST",    $null,  "0.00248362891109477",   "TU_PET*"),
    @(17, "ACT_COST",    $null,  "0.00538337115666179",   "TU_LPG*"),
    @(18, "ACT_COST",    $null,  "0.00252841765861999",   "TU_DSL*"),
    @(19, "ACT_COST",    $null,  "0.00261893423255859",   "TU_FOL*"),
    @(20, "ACT_COST",    $null,  "0.00252841765861999",   "TU_DID*"),
    @(21, "ACT_COST",    $null,  "0.00252841765861999",   "TU_DIJ*"),
    @(22, "ACT_COST",    $null,  "0.00251030880245087",   "TU_JET*"),
    @(23, "ACT_COST",    $null,  "0.0025",                "TU_OTH*"),
    @(24, "ACT_COST",    $null,  "0.000188040616773223",  "TU_COA*"),
    @(25, "ACT_COST",    $null,  "0.000327439423706614",  "TU_COL*")
)

foreach ($row in $rows) {
    $r      = $row[0]
    $attr   = $row[1]
    $year   = $row[2]
    $region = $row[3]
    $psetPn = $row[4]

    $ws.Range("C$r").Value = $attr
    Set-Txt $ws "F$r" $region
    $ws.Range("H$r").Value = $psetPn

    if ($year -ne $null) {
        Set-Txt $ws "D$r" $year
    } elseif ($r -eq 23) {
        # The NCAP_PASTI "Year" value (2010) moved out of row 23 into row 10;
        # row 23's Year cell goes back to being blank, like the other rows
        # in its block.
        $ws.Range("D23").ClearContents()
    }
}
